$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.761.47'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.349.92'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.06'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '650.25'
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.54'
$ws.Range('E7').Value = '  +8.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.463'
$ws.Range('E8').Value = '  +17.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.09'
$ws.Range('E9').Value = '  +23.19%  '
$ws.Range('D11').Value = '3.347.34'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('E12').Value = '  +4.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '43.46'
$ws.Range('E13').Value = '  +20.05%  '
$ws.Range('E14').Value = '  +8.24%  '
$ws.Range('D15').Value = '99.591.65'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '3.985.56'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').Value = '3.362.85'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.46'
$ws.Range('E19').Value = '  +20.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.93'
$ws.Range('E20').Value = '  +10.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '538.35'
$ws.Range('E21').Value = '  +8.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.57'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.31'
$ws.Range('E23').Value = '  +10.41%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.437'
$ws.Range('E25').Value = '  +51.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '103.66'
$ws.Range('E26').Value = '  +15.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.25'
$ws.Range('E27').Value = '  +8.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.67'
$ws.Range('E28').Value = '  +4.95%  '
$ws.Range('D29').Value = '3.528.15'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.150'
$ws.Range('E30').Value = '  +9.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.90'
$ws.Range('E32').Value = '  +12.90%  '
$ws.Range('E33').Value = '  -3.56%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.20'
$ws.Range('E35').Value = '  +4.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.538'
$ws.Range('E36').Value = '  +15.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.77'
$ws.Range('E37').Value = '  +5.57%  '
$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.08'
$ws.Range('E38').Value = '  +4.77%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.156'
$ws.Range('E39').Value = '  +2.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '517.91'
$ws.Range('E40').Value = '  +3.35%  '
$ws.Range('E41').Value = '  -0.47%  '
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.80'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.37'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.817'
$ws.Range('E45').Value = '  +3.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0403'
$ws.Range('E46').Value = '  +22.73%  '
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('E48').Value = '  +3.31%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.79'
$ws.Range('E49').Value = '  +17.61%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '164.21'
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.03'
$ws.Range('E51').Value = '  +7.42%  '
